# Add the "Brent, 1997" study (two comparator arms) as new rows 6 and 7
# at the bottom of the Zhou NMA studies table, then leave a blank
# (but touched/formatted) row 8 below it, matching the author's manual
# data-entry workflow: identity/description columns first for both new
# rows, then the outcome-statistics columns, then the flag/country
# columns, then the age columns, and finally the id column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Brent, 1997 - arm 1 (individual nondirective supportive therapy vs cbt) ---
# --- Row 7: Brent, 1997 - arm 2 (systemic behavior family therapy vs individual NST) ---

# study / description columns (B-J), row 6 then row 7
$ws.Range("B6").Value = "Brent, 1997"
$ws.Range("C6").Value = 1997
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "cbt"
$ws.Range("F6").Value = "cau"
$ws.Range("G6").Value = "cbt"
$ws.Range("H6").Value = "individual nondirective supportive therapy (NST)"
$ws.Range("I6").Value = "bdi"
$ws.Range("J6").Value = "self-report"

$ws.Range("B7").Value = "Brent, 1997"
$ws.Range("C7").Value = 1997
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "other psy"
$ws.Range("F7").Value = "cau"
$ws.Range("G7").Value = "systemic behavior family therapy"
$ws.Range("H7").Value = "individual nondirective supportive therapy (NST)"
$ws.Range("I7").Value = "bdi"
$ws.Range("J7").Value = "self-report"

# post / baseline outcome columns (K-W), row 6 then row 7
$ws.Range("K6").Value = 5.7
$ws.Range("L6").Value = 8.6
$ws.Range("M6").Value = 35
$ws.Range("N6").Value = "NA"
$ws.Range("O6").Value = "NA"
$ws.Range("P6").Value = 24
$ws.Range("Q6").Value = 24.3
$ws.Range("R6").Value = 8.1
$ws.Range("S6").Value = 37
$ws.Range("T6").Value = 25.7
$ws.Range("U6").Value = 7.8
$ws.Range("V6").Value = 35
$ws.Range("W6").Value = 51.9

$ws.Range("K7").Value = 9.1
$ws.Range("L7").Value = 9.1
$ws.Range("M7").Value = 29
$ws.Range("N7").Value = "NA"
$ws.Range("O7").Value = "NA"
$ws.Range("P7").Value = 24
$ws.Range("Q7").Value = 22.6
$ws.Range("R7").Value = 8.2
$ws.Range("S7").Value = 35
$ws.Range("T7").Value = 25.7
$ws.Range("U7").Value = 7.8
$ws.Range("V7").Value = 35
$ws.Range("W7").Value = 78.6

# responder / country / comorbidity / diagnosis columns (X-AE), row 6 then row 7
$ws.Range("X6").Value = "NA"
$ws.Range("Y6").Value = "NA"
$ws.Range("Z6").Value = "NA"
$ws.Range("AA6").Value = "NA"
$ws.Range("AB6").Value = "us"
$ws.Range("AC6").Value = "n"
$ws.Range("AD6").Value = "NA"
$ws.Range("AE6").Value = "mdd"

$ws.Range("X7").Value = "NA"
$ws.Range("Y7").Value = "NA"
$ws.Range("Z7").Value = "NA"
$ws.Range("AA7").Value = "NA"
$ws.Range("AB7").Value = "us"
$ws.Range("AC7").Value = "n"
$ws.Range("AD7").Value = "NA"
$ws.Range("AE7").Value = "mdd"

# age columns (AF-AK), row 6 then row 7
$ws.Range("AF6").Value = 15.7
$ws.Range("AG6").Value = 1.3
$ws.Range("AH6").Value = 15.7
$ws.Range("AI6").Value = 1.5
$ws.Range("AJ6").Value = 15.6
$ws.Range("AK6").Value = "NA"

$ws.Range("AF7").Value = 15.4
$ws.Range("AG7").Value = 1.4
$ws.Range("AH7").Value = 15.7
$ws.Range("AI7").Value = 1.5
$ws.Range("AJ7").Value = 15.6
$ws.Range("AK7").Value = "NA"

# record id column (A), filled in last
$ws.Range("A6").Value = 128
$ws.Range("A7").Value = 129

# Leave a blank spacer row underneath, carrying over the row-border
# formatting from the table above it (A8:D8), same as the source file.
$ws.Range("A8:D8").Style = $ws.Range("X2:AA2").Style

# Match the active selection left behind by the editor after pasting
# the new study's post-control columns.
$ws.Range("N6:O7").Select()
